# "backed out Traffic logic" - remove the Traffic-intent logic and instead
# just list each account's website as a plain hyperlink beneath the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Hyperlinks.Add($ws.Range("B4"), "http://cloudflare.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "http://apple.com") | Out-Null

$ws.Range("B6").Select() | Out-Null
